$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") per row, only where they change
$updates = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 3
    6  = 0
    8  = 0
    9  = 2
    10 = 3
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 2
    21 = 1
    22 = 1
    23 = 2
    24 = 3
    26 = 1
    27 = 2
    28 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
